$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Lugo / Almeria rows (row 47 and row 48): city name and the
# "Casos activos" (column C) value move together so the correct value
# stays attached to the correct province; columns B, D, E are identical
# between the two rows and are left untouched.
$ws.Range("A47").Value = "Almeria"
$ws.Range("C47").Value = 72

$ws.Range("A48").Value = "Lugo"
$ws.Range("C48").Value = 5

# Update the "last updated" timestamp footer in A1.
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 00:16"
